# -----------------------------------------------------------------------
# "run prepare & render with final data"
#
# The upstream data-prep step now includes Russia as an additional
# country column (inserted after Japan, before Saudi Arabia) and the
# policy rows are re-rendered sorted by the recomputed "All" mean
# (column B), descending. This script writes the final rendered
# table - header, row labels, and the full data grid - into the sheet.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Russia" column is inserted into the header (B1:N1); A1 (blank
# corner cell) is left untouched.
$headerLabels = @(
'$ bold(''All'')',  '$ bold(''Europe'')',  'France',  'Germany',  'Italy',  'Poland',  'Spain',  'United Kingdom',  'Switzerland',  'Japan',  'Russia',  'Saudi Arabia',  'USA'
)
for ($i = 0; $i -lt $headerLabels.Length; $i++) {
    $ws.Cells.Item(1, $i + 2).Value2 = $headerLabels[$i]
}

# Final policy-row order (A2:A11), sorted by the recomputed "All" mean
$rowLabels = @(
    'Minimum tax of 2% on billionaires''
wealth, in voluntary countries',
    'Bridgetown initiative: MDBs expanding sustainable
investments in LICs, and at lower interest rates',
    'L&D: Developed countries financing a fund to help
vulnerable countries cope with climate Loss and damage',
    'Debt relief for vulnerable countries, suspending
payments until they are more able to repay',
    'Expand Security Council to new permanent members (e.g.
India, Brazil, African Union), restrict veto use',
    'Raise global minimum tax on profit from 15% to 35%,
allocating revenues to countries based on sales',
    'At least 0.7% of developed countries'' GDP in foreign aid',
    'International levy on shipping carbon emissions,
returned to countries based on population',
    'NCQG: Developing countries providing $300 bn a
year in climate finance for developing countries',
    'International levy on aviation carbon emissions, raising
prices by 30%, returned to countries based on population'
)
for ($r = 0; $r -lt $rowLabels.Length; $r++) {
    $ws.Cells.Item($r + 2, 1).Value2 = $rowLabels[$r]
}

# Final rendered data grid (B2:N11): All, Europe, per-country values,
# including the new Russia column and the recalculated means.
$data = @(
    @(0.686465033529492, 0.829105003523608, 0.906290045675121, 0.837084307502415, 0.926149166801291, 0.58672106292095, 0.872587753455882, 0.878110678654192, 0.666831860457191, 0.43467671690192, 0.663985194495835, 0.703761021752447, 0.640920100928406),
    @(0.501984398521965, 0.575757575757576, 0.460271257364679, 0.530471978022121, 0.814055061924697, 0.281839074996125, 0.708667471645252, 0.668910723233189, 0.500512431803472, 0.265591092846314, 0.596930680039649, 0.794198117239725, 0.422308809717726),
    @(0.430409196660736, 0.457011980267794, 0.33279515263289, 0.405479942619727, 0.676607281778073, 0.331985231281113, 0.663737354016908, 0.463808700943725, 0.292811789467233, 0.192156345298118, 0.666375959335793, 0.85524931908711, 0.321756048101479),
    @(0.344053647187628, 0.354474982381959, 0.192990020769688, 0.148882177492379, 0.63699224731715, 0.396220150149398, 0.530347555433164, 0.421724564755278, 0.228655539408773, 0.150035589499553, 0.38340162387244, 0.829498637932723, 0.258047211971733),
    @(0.335248256052524, 0.461240310077519, 0.374483973296846, 0.443559469608065, 0.652535397934555, 0.335545645974619, 0.566078766657601, 0.468274278957962, 0.346327380067558, 0.147582101842148, 0.00654400091373214, 0.609308133997403, 0.252087215759207),
    @(0.328589024223347, 0.465116279069767, 0.472989584915836, 0.426758680000507, 0.791797773897916, 0.228384151173381, 0.498426814110758, 0.485870949814896, 0.174811044915976, 0.223123315914851, -0.12388965679271, 0.530833388531733, 0.249506556308181),
    @(0.325441357602299, 0.320295983086681, 0.233034713655876, 0.265811882008148, 0.479713313825981, 0.0762488483083057, 0.614577943816408, 0.329627924956973, 0.215601516824756, 0.062262100614088, 0.630975799362228, 0.751452143604012, 0.260707883435066),
    @(0.316135212809635, 0.405567300916138, 0.506526534184808, 0.313599558455335, 0.546322756695602, 0.150011689134294, 0.515525999865488, 0.429632984549167, 0.351403943487063, 0.0188310156555156, 0.356241425568155, 0.625075550597757, 0.253242228564643),
    @(0.31052415491994, 0.36645525017618, 0.261052125627707, 0.384187374065964, 0.549756293248072, 0.200768841642844, 0.560420781364187, 0.338934301090944, 0.234400834579408, 0.0308357678252117, 0.677911376140655, 0.689558758766624, 0.159422230689794),
    @(0.0619953469276037, 0.0916138125440451, 0.153827555592992, 0.0945014460121481, 0.096465759901817, 0.00325139474591023, 0.133069041485609, 0.120480878724206, 0.009637585468826, -0.0928582986299513, -0.0372928283234603, 0.443518183056329, -0.00371332020131786)
)
for ($r = 0; $r -lt $data.Length; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 2).Value2 = $rowVals[$c]
    }
}
